# Regenerate orders with updated distance/size codes.
# Applies the substring substitutions D80->D86, D64->D69, D51->D55, S30->S31
# to every textual value in the used range of the active sheet (shared
# strings such as Condition, Filename_Left, Filename_Right, Distance, Size).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $nv = $v.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
